$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1953.1333
$ws.Range("I62").Value = 1929.9
$ws.Range("J62").Value = 1999.6
$ws.Range("K62").Value = 1929.9
$ws.Range("L62").Value = 1999.6
$ws.Range("M62").Value = -1305.9
$ws.Range("N62").Value = -3247.6

$ws.Range("H65").Value = 1953.1333
$ws.Range("I65").Value = 1929.9
$ws.Range("J65").Value = 1999.6
$ws.Range("K65").Value = 9649.5
$ws.Range("L65").Value = 9998
$ws.Range("M65").Value = -6529.5
$ws.Range("N65").Value = -16238

$ws.Range("H70").Value = 2500
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -9540

$ws.Range("H73").Value = 2500
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -10872

$ws.Range("H74").Value = 5533.5
$ws.Range("J74").Value = 6349.8335
$ws.Range("L74").Value = 6349.8335
$ws.Range("N74").Value = -8221.833500000001

$ws.Range("H76").Value = 3514.25
$ws.Range("I76").Value = 3184.3635
$ws.Range("J76").Value = 4240
$ws.Range("K76").Value = 3184.3635
$ws.Range("L76").Value = 4240
$ws.Range("M76").Value = -2869.3635
$ws.Range("N76").Value = -4870

$ws.Range("H77").Value = 5533.5
$ws.Range("J77").Value = 6349.8335
$ws.Range("L77").Value = 31749.1675
$ws.Range("N77").Value = -41109.1675

$ws.Range("H79").Value = 3514.25
$ws.Range("I79").Value = 3184.3635
$ws.Range("J79").Value = 4240
$ws.Range("K79").Value = 3184.3635
$ws.Range("L79").Value = 4240
$ws.Range("M79").Value = -2092.3635
$ws.Range("N79").Value = -6424

$ws.Range("H86").Value = 6681.8
$ws.Range("I86").Value = 9269.6
$ws.Range("J86").Value = 5387.9
$ws.Range("K86").Value = 9269.6
$ws.Range("L86").Value = 5387.9
$ws.Range("M86").Value = -8146.6
$ws.Range("N86").Value = -7633.9

$ws.Range("H89").Value = 6681.8
$ws.Range("I89").Value = 9269.6
$ws.Range("J89").Value = 5387.9
$ws.Range("K89").Value = 46348
$ws.Range("L89").Value = 26939.5
$ws.Range("M89").Value = -40732
$ws.Range("N89").Value = -38171.5

$ws.Range("H94").Value = 333339000
$ws.Range("I94").Value = 8500
$ws.Range("K94").Value = 8500
$ws.Range("M94").Value = -8049

$ws.Range("H98").Value = 763.0909
$ws.Range("I98").Value = 729.4
$ws.Range("J98").Value = 1100
$ws.Range("K98").Value = 729.4
$ws.Range("L98").Value = 1100
$ws.Range("M98").Value = 768.6
$ws.Range("N98").Value = -4096

$ws.Range("H122").Value = 763.0909
$ws.Range("I122").Value = 729.4
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 2188.2
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = 261.8000000000002
$ws.Range("N122").Value = -8200

$ws.Range("H132").Value = 8500
$ws.Range("I132").Value = 9412.5
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 28237.5
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -25707.5
$ws.Range("N132").Value = -8660

$ws.Range("H137").Value = 1376.449
$ws.Range("I137").Value = 873.3333
$ws.Range("J137").Value = 1598.4117
$ws.Range("K137").Value = 2619.9999
$ws.Range("L137").Value = 4795.2351
$ws.Range("M137").Value = -69.9998999999998
$ws.Range("N137").Value = -9895.2351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26817.193
$ws.Range("I32").Value = 4936.3477
$ws.Range("K32").Value = 4936.3477
$ws.Range("M32").Value = -4649.3477

$ws.Range("H74").Value = 1274.82
$ws.Range("I74").Value = 888.3929000000001
$ws.Range("J74").Value = 1766.6364
$ws.Range("K74").Value = 888.3929000000001
$ws.Range("L74").Value = 1766.6364
$ws.Range("M74").Value = -14.39290000000005
$ws.Range("N74").Value = -3514.6364

$ws.Range("H77").Value = 1274.82
$ws.Range("I77").Value = 888.3929000000001
$ws.Range("J77").Value = 1766.6364
$ws.Range("K77").Value = 4441.9645
$ws.Range("L77").Value = 8833.182000000001
$ws.Range("M77").Value = -73.96450000000004
$ws.Range("N77").Value = -17569.182

$ws.Range("H132").Value = 2807.4285
$ws.Range("I132").Value = 2813.4075
$ws.Range("J132").Value = 2787.25
$ws.Range("K132").Value = 8440.2225
$ws.Range("L132").Value = 8361.75
$ws.Range("M132").Value = -5910.2225
$ws.Range("N132").Value = -13421.75

$ws.Range("H134").Value = 64554.445
$ws.Range("J134").Value = 64554.445
$ws.Range("L134").Value = 64554.445
$ws.Range("N134").Value = -74694.44500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1475.1666
$ws.Range("I80").Value = 613.5
$ws.Range("K80").Value = 613.5
$ws.Range("M80").Value = 384.5

$ws.Range("H83").Value = 1475.1666
$ws.Range("I83").Value = 613.5
$ws.Range("K83").Value = 3067.5
$ws.Range("M83").Value = 1924.5

$ws.Range("H95").Value = 47312
$ws.Range("J95").Value = 47312
$ws.Range("L95").Value = 47312
$ws.Range("N95").Value = -52804

$ws.Range("H107").Value = 25643720
$ws.Range("I107").Value = 83334990
$ws.Range("J107").Value = 3156.889
$ws.Range("K107").Value = 83334990
$ws.Range("L107").Value = 3156.889
$ws.Range("M107").Value = -83333070
$ws.Range("N107").Value = -6996.889

$ws.Range("H134").Value = 2091.85
$ws.Range("I134").Value = 2049.2354
$ws.Range("K134").Value = 6147.706200000001
$ws.Range("M134").Value = -3612.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2107.946
$ws.Range("I31").Value = 1403.5333
$ws.Range("J31").Value = 2588.2273
$ws.Range("K31").Value = 1403.5333
$ws.Range("L31").Value = 2588.2273
$ws.Range("M31").Value = -1108.5333
$ws.Range("N31").Value = -3178.2273

$ws.Range("H34").Value = 2107.946
$ws.Range("I34").Value = 1403.5333
$ws.Range("J34").Value = 2588.2273
$ws.Range("K34").Value = 1403.5333
$ws.Range("L34").Value = 2588.2273
$ws.Range("M34").Value = -1201.5333
$ws.Range("N34").Value = -2992.2273

$ws.Range("H103").Value = 16266.667
$ws.Range("I103").Value = 9400
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 9400
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -8228
$ws.Range("N103").Value = -32344

$ws.Range("H107").Value = 813.5294
$ws.Range("I107").Value = 1218.3334
$ws.Range("J107").Value = 592.7273
$ws.Range("K107").Value = 1218.3334
$ws.Range("L107").Value = 592.7273
$ws.Range("M107").Value = 701.6666
$ws.Range("N107").Value = -4432.7273

$ws.Range("H122").Value = 2799
$ws.Range("I122").Value = 3447.5
$ws.Range("J122").Value = 1872.5714
$ws.Range("K122").Value = 10342.5
$ws.Range("L122").Value = 5617.7142
$ws.Range("M122").Value = -7892.5
$ws.Range("N122").Value = -10517.7142

$ws.Range("H132").Value = 1949.6666
$ws.Range("I132").Value = 1364.1111
$ws.Range("K132").Value = 4092.3333
$ws.Range("M132").Value = -1562.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 233.5
$ws.Range("I8").Value = 233.5
$ws.Range("K8").Value = 700.5
$ws.Range("M8").Value = -561.5

$ws.Range("H107").Value = 253335.02
$ws.Range("I107").Value = 483.60715
$ws.Range("J107").Value = 515551.28
$ws.Range("K107").Value = 1450.82145
$ws.Range("L107").Value = 1546653.84
$ws.Range("M107").Value = 469.1785500000001
$ws.Range("N107").Value = -1550493.84

$ws.Range("H131").Value = 884.2152
$ws.Range("I131").Value = 557.4
$ws.Range("J131").Value = 906.2973
$ws.Range("K131").Value = 1672.2
$ws.Range("L131").Value = 2718.8919
$ws.Range("M131").Value = 3367.8
$ws.Range("N131").Value = -12798.8919

$ws.Range("H141").Value = 2023.1
$ws.Range("I141").Value = 1297.7646
$ws.Range("J141").Value = 6133.3335
$ws.Range("K141").Value = 3893.2938
$ws.Range("L141").Value = 18400.0005
$ws.Range("M141").Value = 1286.7062
$ws.Range("N141").Value = -28760.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 28000
$ws.Range("J26").Value = 28000
$ws.Range("L26").Value = 28000
$ws.Range("N26").Value = -28560

$ws.Range("H50").Value = 28000
$ws.Range("J50").Value = 28000
$ws.Range("L50").Value = 28000
$ws.Range("N50").Value = -28996

$ws.Range("H132").Value = 2200.05
$ws.Range("I132").Value = 1860.2
$ws.Range("J132").Value = 3219.6
$ws.Range("K132").Value = 5580.6
$ws.Range("L132").Value = 9658.799999999999
$ws.Range("M132").Value = -3050.6
$ws.Range("N132").Value = -14718.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 80745.62
$ws.Range("J2").Value = 2226.1904
$ws.Range("L2").Value = 2226.1904
$ws.Range("N2").Value = -2450.1904

$ws.Range("H16").Value = 40709.16
$ws.Range("I16").Value = 50432.4
$ws.Range("K16").Value = 50432.4
$ws.Range("M16").Value = -50262.4

$ws.Range("H75").Value = 7157
$ws.Range("I75").Value = 7157
$ws.Range("K75").Value = 7157
$ws.Range("M75").Value = -6221

$ws.Range("H78").Value = 7157
$ws.Range("I78").Value = 7157
$ws.Range("K78").Value = 21471
$ws.Range("M78").Value = -16791

$ws.Range("H119").Value = 35660
$ws.Range("J119").Value = 35660
$ws.Range("L119").Value = 35660
$ws.Range("N119").Value = -45336

$ws.Range("H122").Value = 2186.652
$ws.Range("I122").Value = 2236.5264
$ws.Range("J122").Value = 1949.75
$ws.Range("K122").Value = 6709.5792
$ws.Range("L122").Value = 5849.25
$ws.Range("M122").Value = -4259.5792
$ws.Range("N122").Value = -10749.25

$ws.Range("H132").Value = 3886.205
$ws.Range("I132").Value = 4624.1904
$ws.Range("J132").Value = 3025.2222
$ws.Range("K132").Value = 13872.5712
$ws.Range("L132").Value = 9075.6666
$ws.Range("M132").Value = -11342.5712
$ws.Range("N132").Value = -14135.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1002740
$ws.Range("I5").Value = 5350
$ws.Range("J5").Value = 1667666.6
$ws.Range("K5").Value = 5350
$ws.Range("L5").Value = 1667666.6
$ws.Range("M5").Value = -5238
$ws.Range("N5").Value = -1667890.6

$ws.Range("H119").Value = 41445
$ws.Range("J119").Value = 41445
$ws.Range("L119").Value = 41445
$ws.Range("N119").Value = -51121
